# Refresh the "quadratic-svm-score" sheet with the latest prediction
# scores (column B) that were previously a stale placeholder copy
# (all 1s) left over from ful-path.csv.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.017968160991749471
$ws.Range("B3").Value = 0.44087867760118815
$ws.Range("B4").Value = 2.0358756717518114
$ws.Range("B5").Value = -0.016432984452297461
$ws.Range("B6").Value = 0.035082233601670509
